# SRDetails.xlsx update
# - Replace the SR numbers in column A with the new set of SR numbers
# - Shared strings get reshuffled as a natural side effect of the rewrite
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SR1234567"
$ws.Range("A3").Value = "SR2345678"
$ws.Range("A4").Value = "SR3456789"
$ws.Range("A5").Value = "SR4567900"
$ws.Range("A6").Value = "SR5679011"
$ws.Range("A7").Value = "SR6790122"
$ws.Range("A8").Value = "SR7901233"
